# Update column F (dSF) values on the active worksheet to reflect
# repulled data / recalculated mean values per the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -5
    "F3"  = -2
    "F4"  = -1
    "F6"  = -6
    "F7"  = -6
    "F8"  = 4
    "F9"  = -3
    "F10" = -4
    "F13" = 4
    "F17" = 10
    "F19" = -2
    "F20" = -3
    "F21" = 5
    "F23" = -4
    "F24" = 0
    "F25" = 3
    "F27" = -6
    "F28" = 6
    "F32" = 5
    "F34" = -4
    "F35" = -7
    "F36" = -2
    "F37" = 1
    "F40" = -1
    "F41" = -2
    "F42" = -1
    "F43" = 0
    "F44" = 2
    "F45" = -1
    "F46" = 4
    "F53" = 2
    "F54" = -3
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
